$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns stay formatted as text so numeric-looking
# strings (e.g. "229.27", "0.0220") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "38.623.05"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "2.100.03"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "229.27"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "61.41"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "0.0845"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "2.414.95"
$ws.Range("E12").Value = "  +3.71%  "
$ws.Range("D13").Value = "14.79"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "22.44"
$ws.Range("E14").Value = "  +6.76%  "
$ws.Range("D15").Value = "0.782"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "5.47"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("D17").Value = "2.113.53"
$ws.Range("E17").Value = "  +5.64%  "
$ws.Range("D18").Value = "38.566.75"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "70.61"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "226.77"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").Value = "169.93"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("D27").Value = "9.43"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = "19.07"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("E30").Value = "  +8.90%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "2.33"
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("D33").Value = "4.74"
$ws.Range("E33").Value = "  +6.07%  "
$ws.Range("D34").Value = "4.47"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "6.47"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "18.23"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").Value = "1.539.92"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "99.98"
$ws.Range("E42").Value = "  +4.83%  "
$ws.Range("D43").Value = "0.0220"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D45").Value = "0.0912"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "7.54"
$ws.Range("E48").Value = "  +6.19%  "
$ws.Range("E49").Value = "  +3.42%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "2.300.73"
$ws.Range("E51").Value = "  +3.74%  "
